# Add FCT and Throughput result table to Sheet2, and make Sheet2 the active sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Header row
$ws2.Range("A1").Value = "num"
$ws2.Range("B1").Value = "RC"
$ws2.Range("C1").Value = "XRC"
$ws2.Range("D1").Value = "ERD"

# Data rows: message size (num), RC, XRC, ERD
$data = @(
    @(16,    1,  7.08,  7.08),
    @(32,    1,  7.61,  7.61),
    @(64,    2,  8.67,  8.67),
    @(128,   3,  11.09, 11.09),
    @(256,   5,  14.27, 16.64),
    @(512,   7,  24.94, 26.85),
    @(1024,  10, 39.78, 47.29),
    @(2048,  17, 84.33, 94.11),
    @(4096,  28, 92.83, 96.01),
    @(8192,  29, 97.3,  98.95),
    @(16384, 33, 97.06, 98.5),
    @(32768, 38, 96.95, 98.19),
    @(65536, 31, 96.32, 97.79),
    @(131072,29, 97.2,  98.44),
    @(262144,37, 96.94, 97.63),
    @(524288,30, 97.37, 98.37)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Select the full data range on Sheet1 before leaving it
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A3:S6").Select()

# Make Sheet2 the active (selected) sheet/tab and leave the cursor near L30
$ws2.Activate()
$ws2.Range("L30").Select()
